$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "275.72"
    "D3"  = "23.05"
    "D4"  = "6.295"
    "D5"  = "0.06272"
    "D6"  = "3.661"
    "D7"  = "6.677"
    "D8"  = "1.366"
    "D9"  = "0.8315"
    "D10" = "0.01381"
    "D11" = "0.1627"
    "D12" = "0.08388"
    "D13" = "0.03447"
    "D14" = "0.03100"
    "D15" = "0.09309"
    "D16" = "3.848"
    "D17" = "0.001661"
    "D19" = "0.006337"
    "D20" = "0.005695"
    "E20" = "19HotbitTokenHTBWorstin24h"
    "D21" = "0.001092"
    "D22" = "0.0001499"
    "D23" = "3.714"
    "D24" = "2.322"
    "D25" = "0.3344"
    "D26" = "0.1239"
    "D27" = "0.0002678"
    "D40" = "0.04703"
    "D41" = "0.007077"
    "D42" = "0.1165"
    "D43" = "0.003698"
    "D44" = "0.01220"
    "D45" = "0.00006247"
    "D47" = "0.7696"
    "E47" = "46CoinbaseStockTokenCOIN"
    "D48" = "0.02967"
    "D50" = "0.01239"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
